$wb = $excel.ActiveWorkbook

# --- "info" sheet: update names ---
$info = $wb.Worksheets.Item("info")
$info.Range("B1").Value = "Daniel"
$info.Range("C1").Value = "Kim"

# --- "items" sheet: update item names and prices ---
$items = $wb.Worksheets.Item("items")

$items.Range("A1").Value = "백김치"
$items.Range("C1").Value = 8000
$items.Range("E1").Value = 8000

$items.Range("A2").Value = "대패삼겹살"
$items.Range("C2").Value = 2000
$items.Range("E2").Value = 2000

$items.Range("A3").Value = "김치"
$items.Range("C3").Value = 8000
$items.Range("E3").Value = 8000
